$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reassign K1..O1 so "Halqa Patwari" is replaced by "Patwari" ---
$ws.Range("K1").Value = "Tehsildar"
$ws.Range("L1").Value = "Patwari"
$ws.Range("M1").Value = "Medical Officer"
$ws.Range("N1").Value = "Local School Headmaster"
$ws.Range("O1").Value = "Counter signed by DC"

# Columns I and J (Date Of Incident / Date Of Report) switch to a text number format
$ws.Range("I1").NumberFormat = "@"
$ws.Range("J1").NumberFormat = "@"

# --- Row 2 data (column H is written later so the shared-string table order matches the target) ---
$ws.Range("A2").Value = "Name1"
$ws.Range("B2").Value = "Father Name 1"
$ws.Range("C2").Value = 12234123
$ws.Range("D2").Value = "Address of affectee"
$ws.Range("E2").Value = "Peshawar"
$ws.Range("F2").Value = "Some Reason"
$ws.Range("G2").Value = 1000
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2012-03-02"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2015-31-04"
$ws.Range("K2").Value = "yes"
$ws.Range("L2").Value = "yes"
$ws.Range("M2").Value = "yes"
$ws.Range("N2").Value = "yes"
$ws.Range("O2").Value = "yes"

# --- Row 3 data (column H is written later so the shared-string table order matches the target) ---
$ws.Range("A3").Value = "Name2"
$ws.Range("B3").Value = "Father Name 2"
$ws.Range("C3").Value = 12234123
$ws.Range("D3").Value = "Address of affectee"
$ws.Range("E3").Value = "Peshawar"
$ws.Range("F3").Value = "Some Reason"
$ws.Range("G3").Value = 1000
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2012-03-02"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2015-31-04"
$ws.Range("K3").Value = "yes"
$ws.Range("L3").Value = "yes"
$ws.Range("M3").Value = "yes"
$ws.Range("N3").Value = "yes"
$ws.Range("O3").Value = "yes"

# --- Column H ("Case type" = injured) written last so "injured" lands at the end of the shared-string table ---
$ws.Range("H2").Value = "injured"
$ws.Range("H3").Value = "injured"

# --- Column widths: col 12 shrinks, new col 13 gets a width (closest achievable values in this engine, which
#     snaps ColumnWidth to 1/6-character increments) ---
$ws.Columns.Item(12).ColumnWidth = 11.333333333333334
$ws.Columns.Item(13).ColumnWidth = 18.666666666666668

# --- Sheet view: drop the frozen/topLeft scroll position, move the active selection to H4 ---
$ws.Range("H4").Select()
